$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep its text type (Price column holds numeric-looking
    # strings like "1.001" or "97.30") instead of Excel auto-converting the
    # assignment to a Number, then restore the default "Normal" style so no
    # extra style index is left attached to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Rows 31 and 32 swap coin data (PancakeSwap <-> Filecoin) with updated prices
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D31") "4.355"
$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D32") "1.478"
$ws.Range("E32").Value = "  +0.58%  "

# Price / Volume(1h) updates for remaining rows
Set-TextValue $ws.Range("D2") "30.228.13"
$ws.Range("E2").Value = "  +0.20%  "
Set-TextValue $ws.Range("D3") "1.861.48"
$ws.Range("E3").Value = "  -0.10%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue $ws.Range("D5") "236.38"
$ws.Range("E5").Value = "  +1.02%  "
Set-TextValue $ws.Range("D7") "0.4680"
$ws.Range("E7").Value = "  +0.23%  "
Set-TextValue $ws.Range("D8") "0.2860"
$ws.Range("E8").Value = "  +1.14%  "
Set-TextValue $ws.Range("D9") "0.06534"
$ws.Range("E9").Value = "  -0.10%  "
Set-TextValue $ws.Range("D10") "21.88"
$ws.Range("E10").Value = "  +8.29%  "
Set-TextValue $ws.Range("D11") "0.07919"
$ws.Range("E11").Value = "  +1.36%  "
Set-TextValue $ws.Range("D12") "97.30"
$ws.Range("E12").Value = "  +1.36%  "
Set-TextValue $ws.Range("D13") "1.868.20"
$ws.Range("E13").Value = "  +0.32%  "
Set-TextValue $ws.Range("D14") "5.157"
$ws.Range("E14").Value = "  +0.58%  "
Set-TextValue $ws.Range("D15") "0.6812"
$ws.Range("E15").Value = "  +1.60%  "
Set-TextValue $ws.Range("D16") "270.33"
$ws.Range("E16").Value = "  -3.38%  "
Set-TextValue $ws.Range("D17") "30.218.61"
$ws.Range("E17").Value = "  +0.05%  "
Set-TextValue $ws.Range("D18") "13.53"
$ws.Range("E18").Value = "  +6.96%  "
Set-TextValue $ws.Range("D20") "0.000007346"
$ws.Range("E20").Value = "  +1.39%  "
Set-TextValue $ws.Range("D21") "2.113.62"
$ws.Range("E21").Value = "  +0.71%  "
Set-TextValue $ws.Range("D22") "5.327"
$ws.Range("E22").Value = "  -2.36%  "
Set-TextValue $ws.Range("D23") "1.001"
$ws.Range("E23").Value = "  +0.08%  "
Set-TextValue $ws.Range("D24") "6.165"
$ws.Range("E24").Value = "  +0.23%  "
Set-TextValue $ws.Range("D25") "167.65"
$ws.Range("E25").Value = "  +1.40%  "
Set-TextValue $ws.Range("D26") "9.213"
$ws.Range("E26").Value = "  -1.15%  "
Set-TextValue $ws.Range("D27") "18.89"
$ws.Range("E27").Value = "  -0.09%  "
Set-TextValue $ws.Range("D28") "1.953"
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("E29").Value = "  +3.05%  "
Set-TextValue $ws.Range("D30") "0.09855"
$ws.Range("E30").Value = "  +2.79%  "
Set-TextValue $ws.Range("D33") "4.054"
$ws.Range("E33").Value = "  -1.62%  "
Set-TextValue $ws.Range("D34") "0.04703"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").Value = "  -0.25%  "
Set-TextValue $ws.Range("D37") "2.712"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.04%  "
Set-TextValue $ws.Range("D39") "2.624"
$ws.Range("E39").Value = "  +3.74%  "
Set-TextValue $ws.Range("D40") "75.46"
$ws.Range("E40").Value = "  +4.17%  "
Set-TextValue $ws.Range("D41") "6.260"
$ws.Range("E41").Value = "  -0.24%  "
Set-TextValue $ws.Range("D42") "1.941"
$ws.Range("E42").Value = "  +0.57%  "
Set-TextValue $ws.Range("D43") "0.8520"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.10%  "
Set-TextValue $ws.Range("D45") "0.4156"
$ws.Range("E45").Value = "  -0.11%  "
Set-TextValue $ws.Range("D46") "103.21"
$ws.Range("E46").Value = "  +0.05%  "
Set-TextValue $ws.Range("D47") "955.05"
$ws.Range("E47").Value = "  -3.36%  "
Set-TextValue $ws.Range("D48") "7.161"
$ws.Range("E48").Value = "  +0.40%  "
Set-TextValue $ws.Range("D49") "9.261"
$ws.Range("E49").Value = "  +0.60%  "
Set-TextValue $ws.Range("D50") "34.12"
$ws.Range("E50").Value = "  +0.29%  "
Set-TextValue $ws.Range("D51") "0.05651"
$ws.Range("E51").Value = "  +0.35%  "
